$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data for subjects 6-9 (rows 8-11)
# Columns: A=ID, B=Condition, C=Mental Demand, D=Physical Demand, E=Temporal Demand, F=Performance, G=Effort, H=Frustration
$data = @(
    @("G", "N/A", 4, 10, 5, 10, 12),
    @("F", 3, 6, 8, 3, 10, 6),
    @("G", 8, 3, 8, 5, 8, 3),
    @("F", 10, 3, 14, 4, 17, 5)
)

$startRow = 8
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $col = 2 + $c
        $ws.Cells.Item($row, $col).Value = $rowData[$c]
    }
}

# Update the selection to H12
$ws.Range("H12").Select()
